# Day 3 update for Tournament1Teams workbook:
#  - Team "Team 2" (row 8) renamed to "2 Bacon Nuggets"
#  - Substitute rows revealed (previously hidden) with substitute player names
#  - Selection moved to K20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reveal substitute rows and fill in the substitute names
# (order matches the shared-string append order in the target workbook:
#  Yandi, then Yante, then the renamed team)
$ws.Rows.Item(31).Hidden = $False
$ws.Range("C31").Value = "Substitute: Yandi"

$ws.Rows.Item(7).Hidden = $False
$ws.Range("C7").Value = "Substitute: Yante"

$ws.Rows.Item(11).Hidden = $False
$ws.Range("C11").Value = "Substitute: Yante"

# Rename "Team 2" -> "2 Bacon Nuggets"
$ws.Range("B8").Value = "2 Bacon Nuggets"

# Move the active selection to K20
$ws.Range("K20").Select()
